$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$cells = @{
    "H32" = 725.1429000000001
    "I32" = 271.14285
    "J32" = 1179.1428
    "K32" = 271.14285
    "L32" = 1179.1428
    "M32" = 54.85714999999999
    "N32" = -1831.1428
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H38" = 1231.7391
    "J38" = 2278.9167
    "L38" = 6836.750100000001
    "N38" = -7580.750100000001
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H58" = 797.95
    "J58" = 1180.7778
    "L58" = 3542.3334
    "N58" = -3842.3334
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H64" = 3833.9
    "I64" = 3797
    "J64" = 3920
    "K64" = 3797
    "L64" = 3920
    "M64" = -3549
    "N64" = -4416
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H67" = 3833.9
    "I67" = 3797
    "J67" = 3920
    "K67" = 3797
    "L67" = 3920
    "M67" = -2939
    "N67" = -5636
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H70" = 1349.4286
    "I70" = 1647.3334
    "J70" = 1126
    "K70" = 4942.0002
    "L70" = 3378
    "M70" = -4672.0002
    "N70" = -3918
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H73" = 1349.4286
    "I73" = 1647.3334
    "J73" = 1126
    "K73" = 4942.0002
    "L73" = 3378
    "M73" = -4006.0002
    "N73" = -5250
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H88" = 650159.75
    "I88" = 453.36365
    "J88" = 1543506
    "K88" = 453.36365
    "L88" = 1543506
    "M88" = -47.36365000000001
    "N88" = -1544318
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H91" = 650159.75
    "I91" = 453.36365
    "J91" = 1543506
    "K91" = 453.36365
    "L91" = 1543506
    "M91" = 950.63635
    "N91" = -1546314
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H107" = 1918.6296
    "I107" = 1328.8572
    "J107" = 3982.8333
    "K107" = 1328.8572
    "L107" = 3982.8333
    "M107" = 591.1428000000001
    "N107" = -7822.8333
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H116" = 3611.8462
    "I116" = 2990.8333
    "J116" = 4144.143
    "K116" = 2990.8333
    "L116" = 4144.143
    "M116" = 451.1667000000002
    "N116" = -11028.143
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H132" = 6066717.5
    "I132" = 7096079.5
    "K132" = 21288238.5
    "M132" = -21285708.5
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$cells = @{
    "H33" = 50001000
    "J33" = 0
    "L33" = 0
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}
$ws.Range("N33").ClearContents()

$cells = @{
    "H133" = 34133.332
    "J133" = 34133.332
    "L133" = 34133.332
    "N133" = -39193.332
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H139" = 30416.4
    "J139" = 30416.4
    "L139" = 30416.4
    "N139" = -40696.4
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$cells = @{
    "H86" = 2924.2
    "I86" = 2931.1875
    "J86" = 2911.7778
    "K86" = 2931.1875
    "L86" = 2911.7778
    "M86" = -1808.1875
    "N86" = -5157.7778
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H89" = 2924.2
    "I89" = 2931.1875
    "J89" = 2911.7778
    "K89" = 14655.9375
    "L89" = 14558.889
    "M89" = -9039.9375
    "N89" = -25790.889
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$cells = @{
    "H31" = 1342.6227
    "I31" = 1310.5333
    "K31" = 1310.5333
    "M31" = -1015.5333
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H34" = 1342.6227
    "I34" = 1310.5333
    "K34" = 1310.5333
    "M34" = -1108.5333
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H99" = 1577
    "I99" = 1477.7273
    "J99" = 1850
    "K99" = 1477.7273
    "L99" = 1850
    "M99" = 20.27269999999999
    "N99" = -4846
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H126" = 1577
    "I126" = 1477.7273
    "J126" = 1850
    "K126" = 4433.1819
    "L126" = 5550
    "M126" = -1963.1819
    "N126" = -10490
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$cells = @{
    "H14" = 227.125
    "I14" = 227.125
    "K14" = 681.375
    "M14" = -508.375
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H131" = 24393566
    "J131" = 4238.355
    "L131" = 12715.065
    "N131" = -22795.065
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$cells = @{
    "H122" = 3360.375
    "I122" = 3376.1428
    "J122" = 3250
    "K122" = 10128.4284
    "L122" = 9750
    "M122" = -7678.428400000001
    "N122" = -14650
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H132" = 2922.879
    "I132" = 3215.4707
    "J132" = 2612
    "K132" = 9646.4121
    "L132" = 7836
    "M132" = -7116.4121
    "N132" = -12896
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$cells = @{
    "H22" = 1451
    "I22" = 0
    "J22" = 1451
    "K22" = 0
    "L22" = 1451
    "N22" = -2041
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}
$ws.Range("M22").ClearContents()

$cells = @{
    "H27" = 1451
    "I27" = 0
    "J27" = 1451
    "K27" = 0
    "L27" = 1451
    "N27" = -1665
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}
$ws.Range("M27").ClearContents()

$cells = @{
    "H61" = 1018.55
    "I61" = 959.26666
    "J61" = 1196.4
    "K61" = 959.26666
    "L61" = 1196.4
    "M61" = -757.26666
    "N61" = -1600.4
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H64" = 19736
    "J64" = 19736
    "L64" = 19736
    "N64" = -20186
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H67" = 19736
    "J67" = 19736
    "L67" = 19736
    "N67" = -21296
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H68" = 1266
    "I68" = 1260.7693
    "J68" = 1300
    "K68" = 1260.7693
    "L68" = 1300
    "M68" = -511.7692999999999
    "N68" = -2798
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H71" = 1266
    "I71" = 1260.7693
    "J71" = 1300
    "K71" = 6303.8465
    "L71" = 6500
    "M71" = -2559.8465
    "N71" = -13988
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H100" = 1258
    "I100" = 1040
    "J100" = 1476
    "K100" = 1040
    "L100" = 1476
    "M100" = -499
    "N100" = -2558
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H113" = 1018.55
    "I113" = 959.26666
    "J113" = 1196.4
    "K113" = 959.26666
    "L113" = 1196.4
    "M113" = 1210.73334
    "N113" = -5536.4
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H132" = 2874.55
    "I132" = 2563.4285
    "K132" = 7690.2855
    "M132" = -5160.2855
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$cells = @{
    "H63" = 14024
    "J63" = 14748.75
    "L63" = 14748.75
    "N63" = -15996.75
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H66" = 14024
    "J66" = 14748.75
    "L66" = 44246.25
    "N66" = -50486.25
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}

$cells = @{
    "H122" = 10871405
    "I122" = 11906625
    "K122" = 35719875
    "M122" = -35717425
}
foreach ($key in $cells.Keys) {
    $ws.Range($key).Value = $cells[$key]
}
